$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.046151399612427
$ws.Range("B1").Value = 1.961359143257141
$ws.Range("C1").Value = 8.249018669128418
$ws.Range("D1").Value = 1.578741312026978
$ws.Range("E1").Value = 0.7908744215965271
